{"js": "// \"Final Changes by Kordik\"\n//\n// The report's closing paragraph (starting \"The super scalar size and\n// number of functional units...\") is followed by a brand-new paragraph\n// that wraps up the optimal-parameter discussion. The hidden \"_GoBack\"\n// bookmark (Word's \"last edit position\" marker) moves from the middle of\n// the old closing paragraph to the very end of the newly added one.\n\nconst body = context.document.body;\n\n// --- 1. Remove the old \"_GoBack\" bookmark (it sat mid-paragraph, between\n//        \" we can limit powe\" and \"r consumption...\") so it doesn't linger\n//        at its old spot once we relocate it. ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Locate the paragraph that ends in \"...reduce unused resources.\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"reduce unused resources\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!targetParagraph) {\n  throw new Error(\"Could not find the 'reduce unused resources.' paragraph\");\n}\n\n// --- 3. Insert a brand-new paragraph right after it, carrying forward the\n//        same paragraph formatting (first-line indent) automatically. ---\nconst newText =\n  \"Using these optmial values as the starting point, the parameters are \" +\n  \"swept again to find a final set of optimal values.  It was found that \" +\n  \"the optimal settings were: super scalar factor = 10, number of \" +\n  \"reservation stations = 18, number of rename register file entries = 28, \" +\n  \"number of reorder buffer entires =25.  These settings yield 1.6 \" +\n  \"instructions per cycle with zero occurances of the register rename \" +\n  \"file full and 26 occurences of the reservation stations being full.\";\n\ntargetParagraph.insertParagraph(newText, Word.InsertLocation.after);\nawait context.sync();\n\n// --- 4. Re-fetch the paragraph collection (the freshly inserted paragraph's\n//        range bounds are only reliable once read again post-sync) and add\n//        the \"_GoBack\" bookmark collapsed at the very end of its text. ---\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nconst newParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\nconst endOfNewParagraph = newParagraph.getRange(Word.RangeLocation.end);\nendOfNewParagraph.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"Final Changes by Kordik\"\n#\n# The report's closing paragraph (starting \"The super scalar size and\n# number of functional units...\") is followed by a brand-new paragraph\n# that wraps up the optimal-parameter discussion. The hidden \"_GoBack\"\n# bookmark (Word's \"last edit position\" marker) moves from the middle of\n# the old closing paragraph to the very end of the newly added one.\n\n$d = $word.ActiveDocument\n\n# --- 1. Locate the paragraph that ends in \"...reduce unused resources.\" ---\n$targetParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $pText = $d.Paragraphs($i).Range.Text\n    if ($pText -like \"*reduce unused resources.*\") {\n        $targetParaIndex = $i\n        break\n    }\n}\nif ($targetParaIndex -eq -1) {\n    throw \"Could not find the 'reduce unused resources.' paragraph\"\n}\n$targetPara = $d.Paragraphs($targetParaIndex)\n\n# --- 2. Remove the old \"_GoBack\" bookmark (it sat mid-paragraph, between\n#        \" we can limit powe\" and \"r consumption...\") so it doesn't linger\n#        at its old spot once we relocate it. ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 3. Insert a brand-new paragraph right after it, carrying forward the\n#        same paragraph formatting (first-line indent) automatically. ---\n$insertAt = $targetPara.Range.Duplicate\n$insertAt.Collapse(0)   # wdCollapseEnd\n$insertAt.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs($targetParaIndex + 1)\n\n$newText = \"Using these optmial values as the starting point, the parameters are swept again to find a final set of optimal values.  It was found that the optimal settings were: super scalar factor = 10, number of reservation stations = 18, number of rename register file entries = 28, number of reorder buffer entires =25.  These settings yield 1.6 instructions per cycle with zero occurances of the register rename file full and 26 occurences of the reservation stations being full.\"\n\n# Append a one-character placeholder so the bookmark's target position is\n# never the literal last character of the paragraph (a boundary position)\n# while we create it; we strip the placeholder right after. This avoids\n# accidentally anchoring the collapsed range to the wrong spot.\n$newPara.Range.Text = $newText + \"X\"\n\n$bmPos = $newPara.Range.Start + $newText.Length\n$bmRange = $d.Range($bmPos, $bmPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$placeholderRange = $d.Range($bmPos, $bmPos + 1)\n$placeholderRange.Delete()\n"}
